$wb = $excel.ActiveWorkbook
$n = $wb.Names.Item(2)
$orig = $n.RefersTo()
Write-Host "orig:" $orig
$n.RefersTo = "=imgproc.hpp!`$B`$1:`$C`$127"
$v = $n.RefersTo()
Write-Host "new refersto:" $v
